{"js": "const replacements = [\n  [\"20\u00d779=1580\", \"70\u00d766=4620\"],\n  [\"63\u00d733=2079\", \"58\u00d728=1624\"],\n  [\"38\u00d799=3762\", \"32\u00d723=736\"],\n  [\"73\u00d711=803\", \"22\u00d745=990\"],\n  [\"39\u00d734=1326\", \"76\u00d722=1672\"],\n  [\"37\u00d741=1517\", \"74\u00d796=7104\"],\n  [\"52\u00d746=2392\", \"40\u00d777=3080\"],\n  [\"67\u00d727=1809\", \"92\u00d755=5060\"],\n  [\"16\u00d765=1040\", \"49\u00d722=1078\"],\n  [\"14\u00d767=938\", \"90\u00d769=6210\"],\n  [\"18\u00d749=882\", \"39\u00d745=1755\"],\n  [\"56\u00d746=2576\", \"65\u00d793=6045\"],\n  [\"92\u00d741=3772\", \"60\u00d731=1860\"],\n  [\"72\u00d753=3816\", \"56\u00d737=2072\"],\n  [\"20\u00d750=1000\", \"54\u00d744=2376\"],\n  [\"25\u00d733=825\", \"45\u00d748=2160\"],\n  [\"67\u00d750=3350\", \"12\u00d751=612\"],\n  [\"62\u00d786=5332\", \"60\u00d723=1380\"],\n  [\"25\u00d751=1275\", \"37\u00d759=2183\"],\n  [\"84\u00d777=6468\", \"15\u00d755=825\"],\n  [\"66\u00d731=2046\", \"92\u00d752=4784\"],\n  [\"93\u00d711=1023\", \"46\u00d794=4324\"],\n  [\"28\u00d760=1680\", \"57\u00d723=1311\"],\n  [\"44\u00d713=572\", \"42\u00d734=1428\"],\n  [\"48\u00d719=912\", \"90\u00d717=1530\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const searchResults = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  for (const range of searchResults.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"20\u00d779=1580\", \"70\u00d766=4620\"),\n    @(\"63\u00d733=2079\", \"58\u00d728=1624\"),\n    @(\"38\u00d799=3762\", \"32\u00d723=736\"),\n    @(\"73\u00d711=803\", \"22\u00d745=990\"),\n    @(\"39\u00d734=1326\", \"76\u00d722=1672\"),\n    @(\"37\u00d741=1517\", \"74\u00d796=7104\"),\n    @(\"52\u00d746=2392\", \"40\u00d777=3080\"),\n    @(\"67\u00d727=1809\", \"92\u00d755=5060\"),\n    @(\"16\u00d765=1040\", \"49\u00d722=1078\"),\n    @(\"14\u00d767=938\", \"90\u00d769=6210\"),\n    @(\"18\u00d749=882\", \"39\u00d745=1755\"),\n    @(\"56\u00d746=2576\", \"65\u00d793=6045\"),\n    @(\"92\u00d741=3772\", \"60\u00d731=1860\"),\n    @(\"72\u00d753=3816\", \"56\u00d737=2072\"),\n    @(\"20\u00d750=1000\", \"54\u00d744=2376\"),\n    @(\"25\u00d733=825\", \"45\u00d748=2160\"),\n    @(\"67\u00d750=3350\", \"12\u00d751=612\"),\n    @(\"62\u00d786=5332\", \"60\u00d723=1380\"),\n    @(\"25\u00d751=1275\", \"37\u00d759=2183\"),\n    @(\"84\u00d777=6468\", \"15\u00d755=825\"),\n    @(\"66\u00d731=2046\", \"92\u00d752=4784\"),\n    @(\"93\u00d711=1023\", \"46\u00d794=4324\"),\n    @(\"28\u00d760=1680\", \"57\u00d723=1311\"),\n    @(\"44\u00d713=572\", \"42\u00d734=1428\"),\n    @(\"48\u00d719=912\", \"90\u00d717=1530\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $find.Text,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $find.Replacement.Text,\n        2\n    )\n}\n"}
